# Centering activity roster update:
# "Clara" (Friday partner for week 2 / Sep 09-12) moves to week 8
# (Oct 21-24), pairing with Gwen. Week 2's Friday slot becomes blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bring D8 up to the same formatting as the rest of the data rows
# (it was a blank, unstyled cell) before writing its new value.
$ws.Range("C8").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("D8").Value = $ws.Range("D2").Value()
$ws.Range("D2").Value = ""
